$wb = $excel.ActiveWorkbook
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 88ccbe24-...md (row 3) and eac24ba7-...md (row 5)
# shared text: 2016-09-02 06:20:20 -> 2016-09-02 06:21:19
$overview.Range("G3").Value = "2016-09-02 06:21:19"
$overview.Range("G5").Value = "2016-09-02 06:21:19"
$dede.Range("H3").Value = "2016-09-02 06:21:19"
$dede.Range("H5").Value = "2016-09-02 06:21:19"

# Priority column "ht" -> "mt" for 88ccbe24-...md (row 3) and eac24ba7-...md (row 5)
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("E5").Value = "mt"
$dede.Range("E3").Value = "mt"
$dede.Range("E5").Value = "mt"

# zh-cn Correspond Handoff Datetime for 88ccbe24-...md (row 3/row5)
# 2016-09-02 06:20:14 -> 2016-09-02 06:21:14
$zhcn.Range("H3").Value = "2016-09-02 06:21:14"
$zhcn.Range("H5").Value = "2016-09-02 06:21:14"

# zh-cn Correspond Handback DateTime for 88ccbe24-...md (row3/row5)
# 2016-09-02 06:20:40 -> 2016-09-02 06:21:32
$zhcn.Range("K3").Value = "2016-09-02 06:21:32"
$zhcn.Range("K5").Value = "2016-09-02 06:21:32"

# de-de Correspond Handback DateTime for 88ccbe24-...md (row3/row5)
# 2016-09-02 06:20:48 -> 2016-09-02 06:21:39
$dede.Range("K3").Value = "2016-09-02 06:21:39"
$dede.Range("K5").Value = "2016-09-02 06:21:39"
